# This sheet reports ligand/receptor interaction scores recomputed with
# updated TPM values. The underlying data changes:
#   - "FAPs -> Itga9/Plg" rows get refreshed numeric scores
#   - three brand new rows are added for the "MuSCs" sending cluster
# We clear all cell contents first (this also resets the shared-string
# table) and then rewrite every cell so the workbook ends up with exactly
# the target values, row count (7) and dimension (A1:T7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.ClearContents()

# Write header row (row 1)
$headers = @("Sending cluster", "Ligand symbol", "Receptor symbol", "Target cluster", "Ligand-expressing cells", "Ligand detection rate", "Ligand average expression value", "Ligand total expression value", "Ligand derived specificity of average expression value", "Ligand derived specificity of total expression value", "Receptor-expressing cells", "Receptor detection rate", "Receptor average expression value", "Receptor total expression value", "Receptor derived specificity of average expression value", "Receptor derived specificity of total expression value", "Edge average expression weight", "Edge total expression weight", "Edge average expression derived specificity", "Edge total expression derived specificity")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i+1).Value2 = $headers[$i]
}

# Row 2
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("B2").Value2 = "Plg"
$ws.Range("C2").Value2 = "Itga9"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.465759
$ws.Range("H2").Value2 = 1.397277
$ws.Range("I2").Value2 = 0.9660838355812051
$ws.Range("J2").Value2 = 0.9660838355812051
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.8366046666666667
$ws.Range("N2").Value2 = 2.509814
$ws.Range("O2").Value2 = 0.08025679986157715
$ws.Range("P2").Value2 = 0.08025679986157715
$ws.Range("Q2").Value2 = 0.389656152942
$ws.Range("R2").Value2 = 3.506905376478
$ws.Range("S2").Value2 = 0.07753479704174558
$ws.Range("T2").Value2 = 0.07753479704174558

# Row 3
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Plg"
$ws.Range("C3").Value2 = "Itga9"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.465759
$ws.Range("H3").Value2 = 1.397277
$ws.Range("I3").Value2 = 0.9660838355812051
$ws.Range("J3").Value2 = 0.9660838355812051
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 7.939250333333333
$ws.Range("N3").Value2 = 23.817751
$ws.Range("O3").Value2 = 0.7616247559221037
$ws.Range("P3").Value2 = 0.7616247559221038
$ws.Range("Q3").Value2 = 3.697777296003
$ws.Range("R3").Value2 = 33.279995664027
$ws.Range("S3").Value2 = 0.7357933654748251
$ws.Range("T3").Value2 = 0.7357933654748252

# Row 4
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "Plg"
$ws.Range("C4").Value2 = "Itga9"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.465759
$ws.Range("H4").Value2 = 1.397277
$ws.Range("I4").Value2 = 0.9660838355812051
$ws.Range("J4").Value2 = 0.9660838355812051
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 1.648242
$ws.Range("N4").Value2 = 4.944726
$ws.Range("O4").Value2 = 0.1581184442163192
$ws.Range("P4").Value2 = 0.1581184442163192
$ws.Range("Q4").Value2 = 0.767683545678
$ws.Range("R4").Value2 = 6.909151911102001
$ws.Range("S4").Value2 = 0.1527556730646344
$ws.Range("T4").Value2 = 0.1527556730646345

# Row 5
$ws.Range("A5").Value2 = "MuSCs"
$ws.Range("B5").Value2 = "Plg"
$ws.Range("C5").Value2 = "Itga9"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.01635133333333333
$ws.Range("H5").Value2 = 0.049054
$ws.Range("I5").Value2 = 0.03391616441879487
$ws.Range("J5").Value2 = 0.03391616441879487
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.8366046666666667
$ws.Range("N5").Value2 = 2.509814
$ws.Range("O5").Value2 = 0.08025679986157715
$ws.Range("P5").Value2 = 0.08025679986157715
$ws.Range("Q5").Value2 = 0.01367960177288889
$ws.Range("R5").Value2 = 0.123116415956
$ws.Range("S5").Value2 = 0.002722002819831564
$ws.Range("T5").Value2 = 0.002722002819831564

# Row 6
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("B6").Value2 = "Plg"
$ws.Range("C6").Value2 = "Itga9"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.01635133333333333
$ws.Range("H6").Value2 = 0.049054
$ws.Range("I6").Value2 = 0.03391616441879487
$ws.Range("J6").Value2 = 0.03391616441879487
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 7.939250333333333
$ws.Range("N6").Value2 = 23.817751
$ws.Range("O6").Value2 = 0.7616247559221037
$ws.Range("P6").Value2 = 0.7616247559221038
$ws.Range("Q6").Value2 = 0.1298173286171111
$ws.Range("R6").Value2 = 1.168355957554
$ws.Range("S6").Value2 = 0.02583139044727857
$ws.Range("T6").Value2 = 0.02583139044727858

# Row 7
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("B7").Value2 = "Plg"
$ws.Range("C7").Value2 = "Itga9"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.01635133333333333
$ws.Range("H7").Value2 = 0.049054
$ws.Range("I7").Value2 = 0.03391616441879487
$ws.Range("J7").Value2 = 0.03391616441879487
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 1.648242
$ws.Range("N7").Value2 = 4.944726
$ws.Range("O7").Value2 = 0.1581184442163192
$ws.Range("P7").Value2 = 0.1581184442163192
$ws.Range("Q7").Value2 = 0.026950954356
$ws.Range("R7").Value2 = 0.242558589204
$ws.Range("S7").Value2 = 0.005362771151684725
$ws.Range("T7").Value2 = 0.005362771151684726
